$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.403199791908264
$ws.Range("B1").Value = 1.971433162689209
$ws.Range("C1").Value = 2.334450483322144
$ws.Range("D1").Value = 4.812489986419678
$ws.Range("E1").Value = 0.9392110705375671
